$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (the "Förändrad" date column) for all data rows (2 through 89)
# from the old date serial 45206 to the new date serial 45208.
$ws.Range("C2:C89").Value = 45208
